$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.512729
$ws.Range("H2").Value = 7.538187
$ws.Range("I2").Value = 0.02190726325199687
$ws.Range("J2").Value = 0.02190726325199687
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 136.6884486666667
$ws.Range("N2").Value = 410.065346
$ws.Range("O2").Value = 0.7423691870207686
$ws.Range("P2").Value = 0.7423691870207685
$ws.Range("Q2").Value = 343.4610289297446
$ws.Range("R2").Value = 3091.149260367702
$ws.Range("S2").Value = 0.01626327721023487
$ws.Range("T2").Value = 0.01626327721023487

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.512729
$ws.Range("H3").Value = 7.538187
$ws.Range("I3").Value = 0.02190726325199687
$ws.Range("J3").Value = 0.02190726325199687
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.8952453333333334
$ws.Range("N3").Value = 2.685736
$ws.Range("O3").Value = 0.004862170554817893
$ws.Range("P3").Value = 0.004862170554817893
$ws.Range("Q3").Value = 2.249508911181334
$ws.Range("R3").Value = 20.245580200632
$ws.Range("S3").Value = 0.0001065168503205032
$ws.Range("T3").Value = 0.0001065168503205033

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.512729
$ws.Range("H4").Value = 7.538187
$ws.Range("I4").Value = 0.02190726325199687
$ws.Range("J4").Value = 0.02190726325199687
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 41.63761133333333
$ws.Range("N4").Value = 124.912834
$ws.Range("O4").Value = 0.2261381995079395
$ws.Range("P4").Value = 0.2261381995079395
$ws.Range("Q4").Value = 104.6240334879953
$ws.Range("R4").Value = 941.616301391958
$ws.Range("S4").Value = 0.004954069067953018
$ws.Range("T4").Value = 0.004954069067953019

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.512729
$ws.Range("H5").Value = 7.538187
$ws.Range("I5").Value = 0.02190726325199687
$ws.Range("J5").Value = 0.02190726325199687
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.903320333333333
$ws.Range("N5").Value = 14.709961
$ws.Range("O5").Value = 0.02663044291647413
$ws.Range("P5").Value = 0.02663044291647413
$ws.Range("Q5").Value = 12.32071519785633
$ws.Range("R5").Value = 110.886436780707
$ws.Range("S5").Value = 0.000583400123488474
$ws.Range("T5").Value = 0.000583400123488474

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 92.89399466666667
$ws.Range("H6").Value = 278.681984
$ws.Range("I6").Value = 0.8098976036382196
$ws.Range("J6").Value = 0.8098976036382197
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 136.6884486666667
$ws.Range("N6").Value = 410.065346
$ws.Range("O6").Value = 0.7423691870207686
$ws.Range("P6").Value = 0.7423691870207685
$ws.Range("Q6").Value = 12697.53602143627
$ws.Range("R6").Value = 114277.8241929265
$ws.Range("S6").Value = 0.6012430255829738
$ws.Range("T6").Value = 0.6012430255829737

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 92.89399466666667
$ws.Range("H7").Value = 278.681984
$ws.Range("I7").Value = 0.8098976036382196
$ws.Range("J7").Value = 0.8098976036382197
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.8952453333333334
$ws.Range("N7").Value = 2.685736
$ws.Range("O7").Value = 0.004862170554817893
$ws.Range("P7").Value = 0.004862170554817893
$ws.Range("Q7").Value = 83.1629152200249
$ws.Range("R7").Value = 748.4662369802242
$ws.Range("S7").Value = 0.003937860280827325
$ws.Range("T7").Value = 0.003937860280827325

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 92.89399466666667
$ws.Range("H8").Value = 278.681984
$ws.Range("I8").Value = 0.8098976036382196
$ws.Range("J8").Value = 0.8098976036382197
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 41.63761133333333
$ws.Range("N8").Value = 124.912834
$ws.Range("O8").Value = 0.2261381995079395
$ws.Range("P8").Value = 0.2261381995079395
$ws.Range("Q8").Value = 3867.884045131406
$ws.Range("R8").Value = 34810.95640618266
$ws.Range("S8").Value = 0.1831487858725418
$ws.Range("T8").Value = 0.1831487858725418

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 92.89399466666667
$ws.Range("H9").Value = 278.681984
$ws.Range("I9").Value = 0.8098976036382196
$ws.Range("J9").Value = 0.8098976036382197
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.903320333333333
$ws.Range("N9").Value = 14.709961
$ws.Range("O9").Value = 0.02663044291647413
$ws.Range("P9").Value = 0.02663044291647413
$ws.Range("Q9").Value = 455.4890128936249
$ws.Range("R9").Value = 4099.401116042624
$ws.Range("S9").Value = 0.0215679319018768
$ws.Range("T9").Value = 0.0215679319018768

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.610639333333333
$ws.Range("H10").Value = 4.831918
$ws.Range("I10").Value = 0.0140423817607685
$ws.Range("J10").Value = 0.0140423817607685
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 136.6884486666667
$ws.Range("N10").Value = 410.065346
$ws.Range("O10").Value = 0.7423691870207686
$ws.Range("P10").Value = 0.7423691870207685
$ws.Range("Q10").Value = 220.1557918348475
$ws.Range("R10").Value = 1981.402126513628
$ws.Range("S10").Value = 0.01042463153157698
$ws.Range("T10").Value = 0.01042463153157698

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.610639333333333
$ws.Range("H11").Value = 4.831918
$ws.Range("I11").Value = 0.0140423817607685
$ws.Range("J11").Value = 0.0140423817607685
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.8952453333333334
$ws.Range("N11").Value = 2.685736
$ws.Range("O11").Value = 0.004862170554817893
$ws.Range("P11").Value = 0.004862170554817893
$ws.Range("Q11").Value = 1.441917346849778
$ws.Range("R11").Value = 12.977256121648
$ws.Range("S11").Value = 0.00006827645511672044
$ws.Range("T11").Value = 0.00006827645511672044

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.610639333333333
$ws.Range("H12").Value = 4.831918
$ws.Range("I12").Value = 0.0140423817607685
$ws.Range("J12").Value = 0.0140423817607685
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 41.63761133333333
$ws.Range("N12").Value = 124.912834
$ws.Range("O12").Value = 0.2261381995079395
$ws.Range("P12").Value = 0.2261381995079395
$ws.Range("Q12").Value = 67.06317455951245
$ws.Range("R12").Value = 603.568571035612
$ws.Range("S12").Value = 0.003175518928183317
$ws.Range("T12").Value = 0.003175518928183318

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.610639333333333
$ws.Range("H13").Value = 4.831918
$ws.Range("I13").Value = 0.0140423817607685
$ws.Range("J13").Value = 0.0140423817607685
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.903320333333333
$ws.Range("N13").Value = 14.709961
$ws.Range("O13").Value = 0.02663044291647413
$ws.Range("P13").Value = 0.02663044291647413
$ws.Range("Q13").Value = 7.897480592799778
$ws.Range("R13").Value = 71.077325335198
$ws.Range("S13").Value = 0.000373954845891483
$ws.Range("T13").Value = 0.000373954845891483

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 17.68108066666667
$ws.Range("H14").Value = 53.04324200000001
$ws.Range("I14").Value = 0.154152751349015
$ws.Range("J14").Value = 0.154152751349015
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 136.6884486666667
$ws.Range("N14").Value = 410.065346
$ws.Range("O14").Value = 0.7423691870207686
$ws.Range("P14").Value = 0.7423691870207685
$ws.Range("Q14").Value = 2416.799487076859
$ws.Range("R14").Value = 21751.19538369173
$ws.Range("S14").Value = 0.1144382526959829
$ws.Range("T14").Value = 0.1144382526959829

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 17.68108066666667
$ws.Range("H15").Value = 53.04324200000001
$ws.Range("I15").Value = 0.154152751349015
$ws.Range("J15").Value = 0.154152751349015
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.8952453333333334
$ws.Range("N15").Value = 2.685736
$ws.Range("O15").Value = 0.004862170554817893
$ws.Range("P15").Value = 0.004862170554817893
$ws.Range("Q15").Value = 15.82890495512356
$ws.Range("R15").Value = 142.460144596112
$ws.Range("S15").Value = 0.0007495169685533449
$ws.Range("T15").Value = 0.0007495169685533449

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 17.68108066666667
$ws.Range("H16").Value = 53.04324200000001
$ws.Range("I16").Value = 0.154152751349015
$ws.Range("J16").Value = 0.154152751349015
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 41.63761133333333
$ws.Range("N16").Value = 124.912834
$ws.Range("O16").Value = 0.2261381995079395
$ws.Range("P16").Value = 0.2261381995079395
$ws.Range("Q16").Value = 736.197964751981
$ws.Range("R16").Value = 6625.781682767829
$ws.Range("S16").Value = 0.03485982563926134
$ws.Range("T16").Value = 0.03485982563926134

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 17.68108066666667
$ws.Range("H17").Value = 53.04324200000001
$ws.Range("I17").Value = 0.154152751349015
$ws.Range("J17").Value = 0.154152751349015
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.903320333333333
$ws.Range("N17").Value = 14.709961
$ws.Range("O17").Value = 0.02663044291647413
$ws.Range("P17").Value = 0.02663044291647413
$ws.Range("Q17").Value = 86.69600234817356
$ws.Range("R17").Value = 780.2640211335621
$ws.Range("S17").Value = 0.004105156045217374
$ws.Range("T17").Value = 0.004105156045217374
